$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 371; everything from row 371 down shifts to 372+
$ws.Rows.Item(371).Insert()

# Populate the new row 371 with the new record
$ws.Cells.Item(371, 1).Value = 5
$ws.Cells.Item(371, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(371, 3).Value = "Maule"
$ws.Cells.Item(371, 4).Value = (Get-Date -Year 2023 -Month 7 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(371, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(371, 5).Value = 7
$ws.Cells.Item(371, 6).Value = "Fruta"
$ws.Cells.Item(371, 7).Value = 100108
$ws.Cells.Item(371, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(371, 9).Value = 100108005
$ws.Cells.Item(371, 10).Value = "Piña"
$ws.Cells.Item(371, 11).Value = "Caramelo"
$ws.Cells.Item(371, 12).Value = "Segunda"
$ws.Cells.Item(371, 13).Value = 150
$ws.Cells.Item(371, 14).Value = 20000
$ws.Cells.Item(371, 15).Value = 20000
$ws.Cells.Item(371, 16).Value = 20000
$ws.Cells.Item(371, 17).Value = '$/caja 14 unidades'
$ws.Cells.Item(371, 18).Value = "Ecuador"
$ws.Cells.Item(371, 19).Value = 1429
$ws.Cells.Item(371, 20).Value = 14
